# Append the next KGID/MRNUMBER record (row 4) to Sheet1 and move the
# selection down to A5, mirroring the cursor position after data entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "KG0002"
$ws.Range("B4").Value = 2

$ws.Range("A5").Select()
